$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 75 (pushes current rows 75-88 down to 77-90),
# mirroring the canonical diff which prepends a new weekly Primera/Segunda
# "Perejil" observation (fecha 45106) ahead of the existing block.
$ws.Rows.Item(75).Insert()
$ws.Rows.Item(75).Insert()

# New row 75: Primera
$ws.Cells.Item(75, 1).Value = 7
$ws.Cells.Item(75, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(75, 3).Value = "Ñuble"
$ws.Cells.Item(75, 4).Value = 45106
$ws.Cells.Item(75, 5).Value = 16
$ws.Cells.Item(75, 6).Value = 100112044
$ws.Cells.Item(75, 7).Value = "Perejil"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 200
$ws.Cells.Item(75, 11).Value = 1500
$ws.Cells.Item(75, 12).Value = 1500
$ws.Cells.Item(75, 13).Value = 1500
$ws.Cells.Item(75, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(75, 15).Value = "Región del Maule"
$ws.Cells.Item(75, 16).Value = 1500
$ws.Cells.Item(75, 17).Value = 1
$ws.Cells.Item(75, 18).Value = "Hortaliza"

# New row 76: Segunda
$ws.Cells.Item(76, 1).Value = 7
$ws.Cells.Item(76, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(76, 3).Value = "Ñuble"
$ws.Cells.Item(76, 4).Value = 45106
$ws.Cells.Item(76, 5).Value = 16
$ws.Cells.Item(76, 6).Value = 100112044
$ws.Cells.Item(76, 7).Value = "Perejil"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Segunda"
$ws.Cells.Item(76, 10).Value = 200
$ws.Cells.Item(76, 11).Value = 1000
$ws.Cells.Item(76, 12).Value = 1000
$ws.Cells.Item(76, 13).Value = 1000
$ws.Cells.Item(76, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(76, 15).Value = "Región del Maule"
$ws.Cells.Item(76, 16).Value = 1000
$ws.Cells.Item(76, 17).Value = 1
$ws.Cells.Item(76, 18).Value = "Hortaliza"
